# Resume content update — skills bullet list rewritten, and the
# experience section's "Undergraduate Research Assistant" bullets plus
# the second entry (Executive Team Member -> Summer Intern @ Viatris
# Egypt) rewritten, per the target diff.

$d = $word.ActiveDocument

# Use Find to locate the old text, then assign Range.Text directly
# (rather than passing the replacement string to Find.Execute) so Word's
# "smart quotes" AutoFormat doesn't mangle straight apostrophes in the
# replacement text, and the run keeps its original formatting (rPr).
function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Text = $old
    $r.Find.MatchWholeWord = $false
    $r.Find.MatchCase = $true
    $found = $r.Find.Execute()
    if (-not $found) {
        throw "Text not found: $old"
    }
    $r.Text = $new
}

# --- Skills bullet list ---
Replace-Text "• Programming languages: Java, Python, JavaScript, C, C++, HTML, CSS, React, Node.JS, Express.JS, SQL" "• Cloud computing and AI integration"
Replace-Text "• Software tools: VS Code, Git, Github, Gitlab, Unity, Unreal Engine" "• Data-driven decision making"
Replace-Text "• Cloud and Data Management: AWS, Azure, GCP basics" "• Cross-functional team collaboration"
Replace-Text "• AI/ML: Familiar with basic AI/ML concepts and integration" "• Adapting to new technologies and tools"
Replace-Text "• Algorithm and Data Structures: Experienced with university-level study" "• Detail-oriented analysis and problem solving"
Replace-Text "• Professional Skills: Adaptable, Excellent communication, Detail-oriented, Leadership, Time Management" "• Strong communication and interpersonal skills"

# --- Undergraduate Research Assistant entry ---
Replace-Text "❖ Undergraduate Research Assistant (Node, React, JS)" "❖ Undergraduate Research Assistant"
Replace-Text "University of Calgary, Calgary, AB" "University of Calgary – Calgary, AB"
Replace-Text "• Developed an automated workflow using Node and React for extracting detailed data insights in a timely manner." "• Developed automated workflows combining cutting-edge technologies like Node, React, and OpenAI's Whisper to enhance data processing efficiency."
Replace-Text "• Collected and processed multi-modal data (videos, spoken recordings, biometric data) for research in information needs." "• Collaborated cross-functionally to create data visualizations for large datasets, honing skills in data-driven insights."
Replace-Text "• Adapted quickly to new tools and technologies to enhance research data analysis processes, showcasing adaptability." "• Demonstrated adaptability by mastering new technology stacks and applied problem-solving to automate data analysis workflows."

# --- Second experience entry: Executive Team Member -> Summer Intern ---
Replace-Text "❖ Executive Team Member" "❖ Summer Intern – Sharpen Up Internship Program (Rotational)"
Replace-Text "Dec 2021 - Apr 2022" "Jun 2025 – Aug 2025"
Replace-Text "Model United Nations (MUN) at Dar Jana International School" "Viatris Egypt – Cairo, Egypt"
Replace-Text "• Organized and prepared event documents, ensuring seamless execution of MUN conferences." "• Gained a broad understanding of pharmaceutical operations by rotating across various departments, contributing to cross-functional initiatives."
Replace-Text "• Managed participant engagement and facilitated communication, enhancing collaborative problem-solving skills." "• Assisted in developing dashboards using Excel and SQL for real-time supply chain performance monitoring."
Replace-Text "• Demonstrated leadership by acting as a spokesman, guiding event procedures effectively." "• Learned and applied pharmaceutical compliance and process improvement principles, enhancing analytical and organizational skills."
